# Generate Report for Handback
# Updates the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps for the second data row (8722a735-...) on the zh-cn and de-de
# sheets, plus the corresponding "Latest HO Xliff Generate Date" roll-up on
# the Overview sheet.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: row 3 is the 8722a735-... file; bump its "Latest HO Xliff
# Generate Date" to the new de-de handback timestamp.
$wsOverview.Range("G3").Value = "2016-10-17 14:36:34"

# zh-cn sheet: row 3 (8722a735-...) handoff/handback timestamps.
$wsZhCn.Range("H3").Value = "2016-10-17 14:36:11"
$wsZhCn.Range("K3").Value = "2016-10-17 14:37:19"

# de-de sheet: row 3 (8722a735-...) handoff/handback timestamps.
$wsDeDe.Range("H3").Value = "2016-10-17 14:36:34"
$wsDeDe.Range("K3").Value = "2016-10-17 14:38:00"
